{"js": "// Replace the 25 three-digit-by-one-digit multiplication problems in the\n// table with their updated values. Each original expression occurs exactly\n// once in the document body, fully inside a single run, so a straightforward\n// search + replace-in-place handles every cell.\nconst replacements = [\n  [\"548\u00d73=1644\", \"189\u00d79=1701\"],\n  [\"439\u00d74=1756\", \"267\u00d74=1068\"],\n  [\"987\u00d77=6909\", \"761\u00d76=4566\"],\n  [\"483\u00d74=1932\", \"973\u00d76=5838\"],\n  [\"486\u00d72=972\", \"283\u00d75=1415\"],\n  [\"878\u00d78=7024\", \"991\u00d79=8919\"],\n  [\"983\u00d72=1966\", \"958\u00d76=5748\"],\n  [\"281\u00d75=1405\", \"698\u00d74=2792\"],\n  [\"671\u00d73=2013\", \"901\u00d72=1802\"],\n  [\"282\u00d72=564\", \"713\u00d72=1426\"],\n  [\"852\u00d74=3408\", \"911\u00d78=7288\"],\n  [\"952\u00d72=1904\", \"729\u00d76=4374\"],\n  [\"597\u00d74=2388\", \"131\u00d72=262\"],\n  [\"977\u00d75=4885\", \"231\u00d76=1386\"],\n  [\"977\u00d73=2931\", \"769\u00d73=2307\"],\n  [\"280\u00d75=1400\", \"463\u00d74=1852\"],\n  [\"709\u00d76=4254\", \"598\u00d75=2990\"],\n  [\"549\u00d78=4392\", \"969\u00d78=7752\"],\n  [\"932\u00d79=8388\", \"444\u00d79=3996\"],\n  [\"685\u00d75=3425\", \"910\u00d75=4550\"],\n  [\"320\u00d79=2880\", \"492\u00d77=3444\"],\n  [\"907\u00d77=6349\", \"429\u00d79=3861\"],\n  [\"637\u00d72=1274\", \"396\u00d73=1188\"],\n  [\"477\u00d73=1431\", \"555\u00d74=2220\"],\n  [\"995\u00d73=2985\", \"295\u00d75=1475\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the 25 three-digit-by-one-digit multiplication problems in the\n# table with their updated values. Each original expression occurs exactly\n# once in the document body, fully inside a single run, so a literal\n# Find/Replace (no wildcards) on the whole document content handles every\n# cell.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"548\u00d73=1644\", \"189\u00d79=1701\"),\n    @(\"439\u00d74=1756\", \"267\u00d74=1068\"),\n    @(\"987\u00d77=6909\", \"761\u00d76=4566\"),\n    @(\"483\u00d74=1932\", \"973\u00d76=5838\"),\n    @(\"486\u00d72=972\", \"283\u00d75=1415\"),\n    @(\"878\u00d78=7024\", \"991\u00d79=8919\"),\n    @(\"983\u00d72=1966\", \"958\u00d76=5748\"),\n    @(\"281\u00d75=1405\", \"698\u00d74=2792\"),\n    @(\"671\u00d73=2013\", \"901\u00d72=1802\"),\n    @(\"282\u00d72=564\", \"713\u00d72=1426\"),\n    @(\"852\u00d74=3408\", \"911\u00d78=7288\"),\n    @(\"952\u00d72=1904\", \"729\u00d76=4374\"),\n    @(\"597\u00d74=2388\", \"131\u00d72=262\"),\n    @(\"977\u00d75=4885\", \"231\u00d76=1386\"),\n    @(\"977\u00d73=2931\", \"769\u00d73=2307\"),\n    @(\"280\u00d75=1400\", \"463\u00d74=1852\"),\n    @(\"709\u00d76=4254\", \"598\u00d75=2990\"),\n    @(\"549\u00d78=4392\", \"969\u00d78=7752\"),\n    @(\"932\u00d79=8388\", \"444\u00d79=3996\"),\n    @(\"685\u00d75=3425\", \"910\u00d75=4550\"),\n    @(\"320\u00d79=2880\", \"492\u00d77=3444\"),\n    @(\"907\u00d77=6349\", \"429\u00d79=3861\"),\n    @(\"637\u00d72=1274\", \"396\u00d73=1188\"),\n    @(\"477\u00d73=1431\", \"555\u00d74=2220\"),\n    @(\"995\u00d73=2985\", \"295\u00d75=1475\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $range = $d.Content\n    $range.Find.ClearFormatting()\n    $range.Find.Replacement.ClearFormatting()\n    $range.Find.Text = $oldText\n    $range.Find.Replacement.Text = $newText\n    $range.Find.MatchWildcards = $false\n    $range.Find.MatchCase = $true\n    $range.Find.Execute($oldText, $true, $true, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
